$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 1.75
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.5
$ws.Range("W5").Value = 4.75
$ws.Range("X5").Value = 6.5
$ws.Range("AH5").Value = 10
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.75
$ws.Range("Y7").Value = 9
$ws.Range("AQ7").Value = 41
$ws.Range("AT7").Value = 2.75
$ws.Range("AU7").Value = 8
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.3
$ws.Range("AD8").Value = 9
$ws.Range("AH8").Value = 23
$ws.Range("AK8").Value = 101
$ws.Range("AU8").Value = 9
$ws.Range("AW8").Value = 9
$ws.Range("AX8").Value = 41
$ws.Range("AY8").Value = 41
$ws.Range("BA8").Value = 151
$ws.Range("L10").Value = 3.1
$ws.Range("N10").Value = 9.5
$ws.Range("BC11").Value = 126
$ws.Range("J12").Value = 3
$ws.Range("L12").Value = 3.5
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.85
$ws.Range("U12").Value = 1.73
$ws.Range("V12").Value = 2
$ws.Range("AC12").Value = 10
$ws.Range("AG12").Value = 201
$ws.Range("AH12").Value = 9.5
$ws.Range("AM12").Value = 29
$ws.Range("AN12").Value = 4.5
$ws.Range("AY12").Value = 23
$ws.Range("BA12").Value = 67
$ws.Range("BB12").Value = 151
$ws.Range("BC12").Value = 151
$ws.Range("G13").Value = 2.9
$ws.Range("I13").Value = 2.4
$ws.Range("J13").Value = 3.5
$ws.Range("W13").Value = 10
$ws.Range("AH13").Value = 8.5
$ws.Range("AM13").Value = 26
$ws.Range("AN13").Value = 5
$ws.Range("M24").Value = 1.05
$ws.Range("N24").Value = 11
$ws.Range("Q24").Value = 1.95
$ws.Range("R24").Value = 1.9
$ws.Range("M25").Value = 1.05
$ws.Range("N25").Value = 11
$ws.Range("Q25").Value = 1.88
$ws.Range("R25").Value = 1.98
$ws.Range("G35").Value = 4.2
$ws.Range("I35").Value = 1.8
$ws.Range("J35").Value = 5.5
$ws.Range("K35").Value = 1.95
$ws.Range("L35").Value = 2.63
$ws.Range("M35").Value = 1.11
$ws.Range("N35").Value = 6.5
$ws.Range("O35").Value = 1.5
$ws.Range("P35").Value = 2.5
$ws.Range("Q35").Value = 2.5
$ws.Range("R35").Value = 1.5
$ws.Range("S35").Value = 1.57
$ws.Range("T35").Value = 2.25
$ws.Range("U35").Value = 2.25
$ws.Range("V35").Value = 1.57
$ws.Range("Y35").Value = 17
$ws.Range("Z35").Value = 51
$ws.Range("AA35").Value = 41
$ws.Range("AB35").Value = 51
$ws.Range("AC35").Value = 6.5
$ws.Range("AD35").Value = 6.5
$ws.Range("AE35").Value = 21
$ws.Range("AF35").Value = 81
$ws.Range("AH35").Value = 5.5
$ws.Range("AI35").Value = 7.5
$ws.Range("AJ35").Value = 9.5
$ws.Range("AK35").Value = 15
$ws.Range("AM35").Value = 41
$ws.Range("AN35").Value = 6
$ws.Range("AO35").Value = 29
$ws.Range("AP35").Value = 41
$ws.Range("AQ35").Value = 101
$ws.Range("AR35").Value = 151
$ws.Range("AT35").Value = 2.25
$ws.Range("AU35").Value = 9.5
$ws.Range("AV35").Value = 81
$ws.Range("AW35").Value = 3.6
$ws.Range("AX35").Value = 11
$ws.Range("AY35").Value = 26
$ws.Range("AZ35").Value = 41
$ws.Range("BA35").Value = 67
$ws.Range("BB35").Value = 251
$ws.Range("G40").Value = 1.73
$ws.Range("H40").Value = 3.5
$ws.Range("I40").Value = 5
$ws.Range("J40").Value = 2.38
$ws.Range("L40").Value = 5.5
$ws.Range("M40").Value = 1.07
$ws.Range("N40").Value = 9
$ws.Range("Q40").Value = 2.1
$ws.Range("R40").Value = 1.7
$ws.Range("U40").Value = 2
$ws.Range("V40").Value = 1.73
$ws.Range("X40").Value = 7.5
$ws.Range("Z40").Value = 13
$ws.Range("AB40").Value = 29
$ws.Range("AE40").Value = 17
$ws.Range("AH40").Value = 12
$ws.Range("AI40").Value = 23
$ws.Range("AJ40").Value = 17
$ws.Range("AN40").Value = 3.6
$ws.Range("AO40").Value = 9
$ws.Range("AS40").Value = 151
$ws.Range("AW40").Value = 6.5
$ws.Range("AY40").Value = 34
$ws.Range("AZ40").Value = 101
$ws.Range("BA40").Value = 126
$ws.Range("BB40").Value = 301
$ws.Range("M42").Value = 1.04
$ws.Range("N42").Value = 12
$ws.Range("O42").Value = 1.22
$ws.Range("P42").Value = 4
$ws.Range("Q42").Value = 1.75
$ws.Range("R42").Value = 2.05
$ws.Range("W42").Value = 6.5
$ws.Range("AA42").Value = 12
$ws.Range("AB42").Value = 34
$ws.Range("AC42").Value = 12
$ws.Range("AE42").Value = 29
$ws.Range("AM42").Value = 67
$ws.Range("AN42").Value = 3.1
$ws.Range("AS42").Value = 151
$ws.Range("BA42").Value = 301
$ws.Range("Q43").Value = 1.75
$ws.Range("L46").Value = 3.5
$ws.Range("AC46").Value = 12
$ws.Range("AL46").Value = 23
$ws.Range("AV46").Value = 51
$ws.Range("AY46").Value = 23
$ws.Range("BB46").Value = 151
